$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("machines")

# Two new machine id rows
$ws.Range("A13").Value = 865006042456120
$ws.Range("A14").Value = 865006042456138

# Format the first new cell (pasted-in look: Arial 12, dark grey #222222, integer format)
$c1 = $ws.Range("A13")
$c1.Font.Size = 12
$c1.Font.Color = 2236962
$c1.Font.Name = "Arial"
$c1.NumberFormat = "0"

# Propagate the same resolved style to the second new cell without re-deriving it
$c1.Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows("13").RowHeight = 15.6
$ws.Rows("14").RowHeight = 15.6

# Column needed to grow to fit the longer values
$ws.Columns("A").ColumnWidth = 25

$ws.Range("A14").Select()
